$d = $word.ActiveDocument

$replacements = @(
    @{old = "484÷5=96, 4";   new = "543÷5=108, 3"},
    @{old = "314÷4=78, 2";   new = "718÷5=143, 3"},
    @{old = "382÷3=127, 1";  new = "828÷4=207, 0"},
    @{old = "254÷2=127, 0";  new = "829÷8=103, 5"},
    @{old = "295÷5=59, 0";   new = "572÷3=190, 2"},
    @{old = "577÷9=64, 1";   new = "919÷3=306, 1"},
    @{old = "980÷8=122, 4";  new = "425÷8=53, 1"},
    @{old = "998÷4=249, 2";  new = "802÷6=133, 4"},
    @{old = "627÷7=89, 4";   new = "338÷6=56, 2"},
    @{old = "927÷3=309, 0";  new = "634÷8=79, 2"},
    @{old = "168÷4=42, 0";   new = "851÷7=121, 4"},
    @{old = "343÷7=49, 0";   new = "658÷2=329, 0"},
    @{old = "696÷7=99, 3";   new = "165÷5=33, 0"},
    @{old = "497÷6=82, 5";   new = "851÷4=212, 3"},
    @{old = "123÷2=61, 1";   new = "109÷7=15, 4"},
    @{old = "829÷5=165, 4";  new = "640÷5=128, 0"},
    @{old = "946÷9=105, 1";  new = "585÷5=117, 0"},
    @{old = "852÷8=106, 4";  new = "510÷2=255, 0"},
    @{old = "779÷6=129, 5";  new = "877÷5=175, 2"},
    @{old = "508÷7=72, 4";   new = "213÷8=26, 5"},
    @{old = "151÷3=50, 1";   new = "538÷3=179, 1"},
    @{old = "829÷7=118, 3";  new = "797÷8=99, 5"},
    @{old = "579÷9=64, 3";   new = "418÷6=69, 4"},
    @{old = "680÷8=85, 0";   new = "239÷7=34, 1"},
    @{old = "546÷8=68, 2";   new = "901÷5=180, 1"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
